# Applies the PAL_Deterministic Schedule_V1_H2 data edit:
#  - Shift every Pre_ISI (column I) value down by 2 seconds.
#    A handful of rows had their Pre_ISI previously clipped at the 5.01s
#    ceiling; those recover their true (un-clipped) underlying value here,
#    so they don't follow the simple "-2" rule like the rest of the column.
#  - Clear the worksheet AutoFilter (and the _FilterDatabase defined name
#    that goes with it).
#  - Update the active selection left behind in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Pre_ISI (column I) was clipped at 5.01 in the prior version;
# these get an explicit corrected value instead of old_value - 2.
$specialPreIsi = @{
    8  = 2.8200000000000003
    12 = 2.96
    14 = 2.7799999999999994
    23 = 2.7300000000000004
    40 = 2.7200000000000006
    43 = 2.76
    47 = 2.7200000000000006
    85 = 2.75
}

for ($row = 2; $row -le 96; $row++) {
    $cell = $ws.Cells.Item($row, 9)   # column I = Pre_ISI
    if ($specialPreIsi.ContainsKey($row)) {
        $cell.Value = $specialPreIsi[$row]
    } else {
        $cell.Value = $cell.Value2 - 2
    }
}

# Remove the AutoFilter (and its backing _FilterDatabase defined name).
$ws.AutoFilterMode = $false
foreach ($n in @($wb.Names)) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.Delete()
    }
}

# Move the remembered selection from S9 to N11.
[void]$ws.Range("N11").Select()
